$d = $word.ActiveDocument

# wdFindContinue = 1
$wdFindContinue = 1

# Use Find to locate text, then assign .Text directly on the matched Range.
# (Find.Execute's built-in "replace" argument runs the text through
#  AutoFormatAsYouType, which mangles straight quotes into curly ones;
#  assigning Range.Text avoids that and keeps formatting/run identity.)
function Replace-Text($findText, $replaceText) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false)
    if (-not $ok) {
        Write-Host "NOT FOUND:" $findText
        return
    }
    $rng.Text = $replaceText
}

# --- Prose / description occurrences (property name used in descriptive sentence) ---
Replace-Text "fwVersion (string)" "FWVersion (string)"
Replace-Text "msSens (" "MotionSensitivity ("
Replace-Text "syncClock (" "ClockSync ("
Replace-Text "maxBrightness (number):" "BrightLevel (number):"
Replace-Text "dimmedBrightness (number):" "DimLevel (number):"
Replace-Text "msBrightness (number):" "MotionLevel (number):"
Replace-Text "holdTime (number):" "HoldTime (number):"
Replace-Text "groupId (number):" "GroupId (number):"
Replace-Text "zoneId (number): " "ZoneId (number): "
Replace-Text "triggerers (array of strings):" "Triggerers (array of strings):"
Replace-Text "triggerees (array of strings): " "Triggerees (array of strings): "

# --- Realign the wrapped continuation line under "MotionLevel (number):" (2 fewer spaces) ---
Replace-Text "                                                           sensor." "                                                         sensor."

# --- JSON sample payload occurrences (quoted key names) ---
Replace-Text '"fwVersion": "1.0",' '"FWVersion": "1.0",'
Replace-Text '"msSens": 3,' '"MotionSensitivity": 3,'
Replace-Text '"syncClock": "Enable",' '"ClockSync": "Enable",'
Replace-Text '"maxBrightness": 100,' '"BrightLevel": 100,'
Replace-Text '"dimmedBrightness": 100,' '"DimLevel": 100,'
Replace-Text '"msBrightness": 100,' '"MotionLevel": 100,'
Replace-Text '"holdTime": 0,' '"HoldTime": 0,'
Replace-Text '"groupId": 255,' '"GroupId": 255,'
Replace-Text '"zoneId": 255,' '"ZoneId": 255,'
Replace-Text '"triggerers": [],' '"Triggerers": [],'
Replace-Text '"triggerees": []' '"Triggerees": []'

Write-Host "Done"
